# ------------------------------------------------------------------
# fix: changed field types and added data input protection for
# absence types
#
#  1) Add a new "Типи відсутності" (absence types) reference sheet
#     listing the allowed values, styled with a bigger serif font.
#  2) Reformat the data-entry columns of the main sheet (A-D, F-P) as
#     Text so identifiers / dates / order numbers are not mangled.
#  3) Fix up a couple of date values on row 2 (now entered as protected
#     text) and drop two stray leftover values.
#  4) Add a dropdown data-validation list on column F restricted to the
#     values from the new reference sheet, with a warning message.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --------------------------------------------------------------
# 1. Create the "Типи відсутності" sheet right after the main sheet
# --------------------------------------------------------------
$typesSheet = $wb.Worksheets.Add([System.Type]::Missing, $ws)
$typesSheet.Name = "Типи відсутності"

$typesSheet.Range("A1").Value = "Відрядження"
$typesSheet.Range("A2").Value = "Відпустка основна"
$typesSheet.Range("A3").Value = "Лікарняний"
$typesSheet.Range("A4").Value = "Стаціонар"
$typesSheet.Range("A5").Value = "Відпустка за сімейними обставинами"
$typesSheet.Range("A6").Value = "Відпустка за станом здоров'я"
$typesSheet.Range("A7").Value = "ВЛК"

# Style the list with a bigger, serif font
$typesListRange = $typesSheet.Range("A1:A7")
$typesListRange.Font.Name = "Times New Roman"
$typesListRange.Font.Size = 14

$typesSheet.Range("E8").Select()

# --------------------------------------------------------------
# 2. Reformat main-sheet data columns (A-D, F-P) as Text.
#    Only cells that already hold data are touched individually so no
#    stray blank formatted cells get materialised.
# --------------------------------------------------------------
$textCells = @("A1", "B1", "C1", "D1", "F1", "G1", "H1", "I1", "J1", "K1", "L1", "M1", "N1", "O1", "P1", "A2", "B2", "C2", "D2", "F2", "G2", "H2", "J2", "A3", "B3", "C3", "D3", "F3", "H3", "J3", "L3", "M3", "O3", "P3", "A4", "B4", "C4", "D4", "F4", "G4", "H4", "J4", "K4", "L4", "P4")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# --------------------------------------------------------------
# 3. Row 2 data fix-up: correct the date values, drop the stray ones
# --------------------------------------------------------------
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "21.12.2024"

$ws.Range("L2").NumberFormat = "@"
$ws.Range("L2").Value = "'04.01.2025"

$ws.Range("K2").ClearContents()
$ws.Range("M2").ClearContents()

# --------------------------------------------------------------
# 4. Data validation: column F restricted to the list on the new sheet
# --------------------------------------------------------------
$val = $ws.Range("F1:F1048576").Validation
$val.Delete()
$val.Add(3, 2, 1, "='Типи відсутності'!`$A`$1:`$A`$7")
$val.ErrorTitle = "Некоректне значення"
$val.ErrorMessage = "Краще використати одне з значень, вказаних на аркуші ""Типи відсутності"""
$val.ShowError = $true
$val.ShowInput = $true

$ws.Select()
$ws.Range("F11").Select()

Write-Host "edit complete"
